# Add a new "ROUND OBS" column (G) next to the existing "ROUND" header (F),
# and rename the existing "ROUND" header to "ROUND MEAN" so the two rounded
# columns are distinguishable (rounded mean vs. rounded observation count).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone F1's formatting (bold header style) onto the new G1 header cell
# before touching any values, so the shared-string bookkeeping below stays
# clean (copy/paste only moves the *format*, not the content).
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)   # xlPasteFormats

# Rename the existing header in place, then add the brand new header text.
# Doing F1 first keeps its shared-string slot stable and appends "ROUND OBS"
# as a fresh entry for G1.
$ws.Range("F1").Value = "ROUND MEAN"
$ws.Range("G1").Value = "ROUND OBS"

# Widen the two rounding columns so the two-word headers are readable.
$ws.Columns("F").ColumnWidth = 14.166666666666666
$ws.Columns("G").ColumnWidth = 13.830729166666666

# The header row now wraps onto a second line ("ROUND OBS"), so bump its
# height to fit two lines of text.
$ws.Rows("1").RowHeight = 28.5

# Leave the selection where the author left it after adding the column.
$ws.Range("G2").Select()
